$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry updates one cell (Price in column D and/or Volume(1h) % in column E)
# for the given row. Values are written with a leading quote so Excel stores them
# as literal text (matching the workbook's existing inline-string cells) instead of
# auto-converting number-looking strings like "63.908.00" or "1.00" into numerics.
# The style is then reset to "Normal" so no stray number-format / quote-prefix
# formatting is left behind on the cell.
$updates = @(
    @{ Addr = "D2"; Value = "63.908.00" },
    @{ Addr = "E2"; Value = "  +5.51%  " },
    @{ Addr = "D3"; Value = "2.721.81" },
    @{ Addr = "E3"; Value = "  +3.70%  " },
    @{ Addr = "E4"; Value = "  +0.03%  " },
    @{ Addr = "D5"; Value = "575.81" },
    @{ Addr = "E5"; Value = "  -0.43%  " },
    @{ Addr = "D6"; Value = "153.67" },
    @{ Addr = "E6"; Value = "  +6.41%  " },
    @{ Addr = "E7"; Value = "  -0.08%  " },
    @{ Addr = "E8"; Value = "  +1.59%  " },
    @{ Addr = "D9"; Value = "2.744.61" },
    @{ Addr = "E9"; Value = "  +4.17%  " },
    @{ Addr = "D10"; Value = "6.69" },
    @{ Addr = "E10"; Value = "  +2.33%  " },
    @{ Addr = "E11"; Value = "  +5.77%  " },
    @{ Addr = "E12"; Value = "  +5.02%  " },
    @{ Addr = "D13"; Value = "0.390" },
    @{ Addr = "E13"; Value = "  +4.14%  " },
    @{ Addr = "D14"; Value = "3.210.12" },
    @{ Addr = "E14"; Value = "  +3.86%  " },
    @{ Addr = "E15"; Value = "  +0.31%  " },
    @{ Addr = "D16"; Value = "63.757.99" },
    @{ Addr = "E16"; Value = "  +5.28%  " },
    @{ Addr = "E17"; Value = "  +6.90%  " },
    @{ Addr = "D18"; Value = "2.746.40" },
    @{ Addr = "E18"; Value = "  +4.47%  " },
    @{ Addr = "D19"; Value = "11.92" },
    @{ Addr = "E19"; Value = "  +3.32%  " },
    @{ Addr = "E20"; Value = "  +3.08%  " },
    @{ Addr = "D21"; Value = "359.85" },
    @{ Addr = "E21"; Value = "  +3.27%  " },
    @{ Addr = "D22"; Value = "6.91" },
    @{ Addr = "E22"; Value = "  +0.53%  " },
    @{ Addr = "D23"; Value = "0.994" },
    @{ Addr = "E23"; Value = "  -0.46%  " },
    @{ Addr = "E24"; Value = "  +0.68%  " },
    @{ Addr = "D25"; Value = "65.96" },
    @{ Addr = "E25"; Value = "  +3.26%  " },
    @{ Addr = "D26"; Value = "0.168" },
    @{ Addr = "E26"; Value = "  +4.63%  " },
    @{ Addr = "E27"; Value = "  +5.70%  " },
    @{ Addr = "D28"; Value = "1.00" },
    @{ Addr = "E28"; Value = "  +0.14%  " },
    @{ Addr = "D29"; Value = "0.0₃0899" },
    @{ Addr = "E29"; Value = "  +12.64%  " },
    @{ Addr = "E30"; Value = "  -1.35%  " },
    @{ Addr = "D31"; Value = "7.10" },
    @{ Addr = "E31"; Value = "  +6.95%  " },
    @{ Addr = "D32"; Value = "171.98" },
    @{ Addr = "E32"; Value = "  +1.79%  " },
    @{ Addr = "E33"; Value = "  +14.22%  " },
    @{ Addr = "D34"; Value = "0.998" },
    @{ Addr = "E34"; Value = "  -0.04%  " },
    @{ Addr = "D35"; Value = "20.47" },
    @{ Addr = "E35"; Value = "  +4.39%  " },
    @{ Addr = "E36"; Value = "  +7.94%  " },
    @{ Addr = "E37"; Value = "  +9.60%  " },
    @{ Addr = "E38"; Value = "  +10.25%  " },
    @{ Addr = "D39"; Value = "0.994" },
    @{ Addr = "E39"; Value = "  +13.80%  " },
    @{ Addr = "D40"; Value = "341.20" },
    @{ Addr = "E40"; Value = "  +3.02%  " },
    @{ Addr = "D41"; Value = "4.20" },
    @{ Addr = "E41"; Value = "  +5.14%  " },
    @{ Addr = "D42"; Value = "39.19" },
    @{ Addr = "E42"; Value = "  +2.34%  " },
    @{ Addr = "E43"; Value = "  +7.00%  " },
    @{ Addr = "D44"; Value = "21.69" },
    @{ Addr = "E44"; Value = "  +7.95%  " },
    @{ Addr = "D45"; Value = "21.73" },
    @{ Addr = "E45"; Value = "  +4.88%  " },
    @{ Addr = "D46"; Value = "0.0589" },
    @{ Addr = "E46"; Value = "  +5.74%  " },
    @{ Addr = "D47"; Value = "0.645" },
    @{ Addr = "E47"; Value = "  +5.85%  " },
    @{ Addr = "D48"; Value = "139.01" },
    @{ Addr = "E48"; Value = "  +4.70%  " },
    @{ Addr = "E49"; Value = "  +4.68%  " },
    @{ Addr = "E50"; Value = "  +1.03%  " },
    @{ Addr = "E51"; Value = "  -0.05%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}
